$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 4499.5
$ws.Range("I86").Value = 4499
$ws.Range("K86").Value = 4499
$ws.Range("M86").Value = -3376
$ws.Range("H89").Value = 4499.5
$ws.Range("I89").Value = 4499
$ws.Range("K89").Value = 22495
$ws.Range("M89").Value = -16879
$ws.Range("H106").Value = 10035.444
$ws.Range("I106").Value = 4576.533
$ws.Range("K106").Value = 4576.533
$ws.Range("M106").Value = -3945.533
$ws.Range("H132").Value = 2253.8704
$ws.Range("I132").Value = 1061.8431
$ws.Range("K132").Value = 3185.5293
$ws.Range("M132").Value = -655.5293000000001
$ws.Range("H137").Value = 3229.88
$ws.Range("I137").Value = 2551.762
$ws.Range("J137").Value = 6790
$ws.Range("K137").Value = 7655.286
$ws.Range("L137").Value = 20370
$ws.Range("M137").Value = -5105.286
$ws.Range("N137").Value = -25470

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 483.2
$ws.Range("I5").Value = 483.2
$ws.Range("K5").Value = 483.2
$ws.Range("M5").Value = -371.2
$ws.Range("H32").Value = 13670.113
$ws.Range("I32").Value = 11940.119
$ws.Range("J32").Value = 50000
$ws.Range("K32").Value = 11940.119
$ws.Range("L32").Value = 50000
$ws.Range("M32").Value = -11653.119
$ws.Range("N32").Value = -50574
$ws.Range("H45").Value = 23097.2
$ws.Range("I45").Value = 35003.668
$ws.Range("K45").Value = 35003.668
$ws.Range("M45").Value = -34626.668
$ws.Range("H63").Value = 3840.4
$ws.Range("I63").Value = 3450.5
$ws.Range("K63").Value = 3450.5
$ws.Range("M63").Value = -2764.5
$ws.Range("H64").Value = 25091
$ws.Range("J64").Value = 25091
$ws.Range("L64").Value = 25091
$ws.Range("N64").Value = -25587
$ws.Range("H66").Value = 3840.4
$ws.Range("I66").Value = 3450.5
$ws.Range("K66").Value = 17252.5
$ws.Range("M66").Value = -13820.5
$ws.Range("H67").Value = 25091
$ws.Range("J67").Value = 25091
$ws.Range("L67").Value = 25091
$ws.Range("N67").Value = -26807
$ws.Range("H74").Value = 2482.0908
$ws.Range("I74").Value = 2225.375
$ws.Range("K74").Value = 2225.375
$ws.Range("M74").Value = -1351.375
$ws.Range("H77").Value = 2482.0908
$ws.Range("I77").Value = 2225.375
$ws.Range("K77").Value = 11126.875
$ws.Range("M77").Value = -6758.875
$ws.Range("H122").Value = 2842.394
$ws.Range("J122").Value = 6000
$ws.Range("L122").Value = 18000
$ws.Range("N122").Value = -22900
$ws.Range("H132").Value = 3136.5715
$ws.Range("I132").Value = 3175.1516
$ws.Range("K132").Value = 9525.4548
$ws.Range("M132").Value = -6995.4548

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 483.2
$ws.Range("I4").Value = 483.2
$ws.Range("K4").Value = 483.2
$ws.Range("M4").Value = -368.2
$ws.Range("H82").Value = 15332.333
$ws.Range("J82").Value = 24998
$ws.Range("L82").Value = 24998
$ws.Range("N82").Value = -25764
$ws.Range("H85").Value = 15332.333
$ws.Range("J85").Value = 24998
$ws.Range("L85").Value = 24998
$ws.Range("N85").Value = -27650
$ws.Range("H108").Value = 80250
$ws.Range("J108").Value = 80250
$ws.Range("L108").Value = 80250
$ws.Range("N108").Value = -87930
$ws.Range("H134").Value = 1935.2727
$ws.Range("I134").Value = 1778.25
$ws.Range("J134").Value = 3505.5
$ws.Range("K134").Value = 5334.75
$ws.Range("L134").Value = 10516.5
$ws.Range("M134").Value = -2799.75
$ws.Range("N134").Value = -15586.5

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 45234.625
$ws.Range("I31").Value = 55164.05
$ws.Range("J31").Value = 7502.8
$ws.Range("K31").Value = 55164.05
$ws.Range("L31").Value = 7502.8
$ws.Range("M31").Value = -54869.05
$ws.Range("N31").Value = -8092.8
$ws.Range("H34").Value = 45234.625
$ws.Range("I34").Value = 55164.05
$ws.Range("J34").Value = 7502.8
$ws.Range("K34").Value = 55164.05
$ws.Range("L34").Value = 7502.8
$ws.Range("M34").Value = -54962.05
$ws.Range("N34").Value = -7906.8
$ws.Range("H51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").ClearContents()
$ws.Range("H60").Value = 0
$ws.Range("J60").Value = 0
$ws.Range("L60").Value = 0
$ws.Range("N60").ClearContents()
$ws.Range("H61").Value = 0
$ws.Range("J61").Value = 0
$ws.Range("L61").Value = 0
$ws.Range("N61").ClearContents()
$ws.Range("H86").Value = 2615.8572
$ws.Range("J86").Value = 3354.6667
$ws.Range("L86").Value = 3354.6667
$ws.Range("N86").Value = -5600.6667
$ws.Range("H89").Value = 2615.8572
$ws.Range("J89").Value = 3354.6667
$ws.Range("L89").Value = 16773.3335
$ws.Range("N89").Value = -28005.3335
$ws.Range("H141").Value = 698553.1
$ws.Range("J141").Value = 780122.5
$ws.Range("L141").Value = 780122.5
$ws.Range("N141").Value = -790482.5

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 9
$ws.Range("I2").Value = 9
$ws.Range("K2").Value = 54
$ws.Range("M2").Value = 59
$ws.Range("H17").Value = 1000
$ws.Range("I17").Value = 1000
$ws.Range("K17").Value = 3000
$ws.Range("M17").Value = -2831
$ws.Range("H34").Value = 145027.72
$ws.Range("J34").Value = 4631.6665
$ws.Range("L34").Value = 13894.9995
$ws.Range("N34").Value = -14062.9995
$ws.Range("H39").Value = 154570.86
$ws.Range("J39").Value = 13666.167
$ws.Range("L39").Value = 40998.501
$ws.Range("N39").Value = -41586.501
$ws.Range("H55").Value = 8406006
$ws.Range("J55").Value = 13897710
$ws.Range("L55").Value = 41693130
$ws.Range("N55").Value = -41693484
$ws.Range("H98").Value = 739.5
$ws.Range("I98").Value = 602
$ws.Range("K98").Value = 1806
$ws.Range("M98").Value = -308

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 41328.266
$ws.Range("J46").Value = 44554.715
$ws.Range("L46").Value = 44554.715
$ws.Range("N46").Value = -44866.715
$ws.Range("H57").Value = 16938.705
$ws.Range("J57").Value = 19279.715
$ws.Range("L57").Value = 19279.715
$ws.Range("N57").Value = -20919.715
$ws.Range("H80").Value = 5666.3335
$ws.Range("I80").Value = 3499.5
$ws.Range("J80").Value = 10000
$ws.Range("K80").Value = 3499.5
$ws.Range("L80").Value = 10000
$ws.Range("M80").Value = -2501.5
$ws.Range("N80").Value = -11996
$ws.Range("H83").Value = 5666.3335
$ws.Range("I83").Value = 3499.5
$ws.Range("J83").Value = 10000
$ws.Range("K83").Value = 17497.5
$ws.Range("L83").Value = 50000
$ws.Range("M83").Value = -12505.5
$ws.Range("N83").Value = -59984
$ws.Range("H102").Value = 1794.9259
$ws.Range("I102").Value = 1702.3182
$ws.Range("K102").Value = 1702.3182
$ws.Range("M102").Value = -80.31819999999993
$ws.Range("H122").Value = 2740.0908
$ws.Range("I122").Value = 2652.4707
$ws.Range("J122").Value = 3038
$ws.Range("K122").Value = 7957.4121
$ws.Range("L122").Value = 9114
$ws.Range("M122").Value = -5507.4121
$ws.Range("N122").Value = -14014
$ws.Range("H132").Value = 9566.902
$ws.Range("I132").Value = 8846.656000000001
$ws.Range("J132").Value = 12127.777
$ws.Range("K132").Value = 26539.968
$ws.Range("L132").Value = 36383.331
$ws.Range("M132").Value = -24009.968
$ws.Range("N132").Value = -41443.331

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 14517.177
$ws.Range("I40").Value = 14924.8125
$ws.Range("K40").Value = 14924.8125
$ws.Range("M40").Value = -14788.8125
$ws.Range("H55").Value = 2479.2856
$ws.Range("I55").Value = 1514.72
$ws.Range("J55").Value = 3897.7646
$ws.Range("K55").Value = 1514.72
$ws.Range("L55").Value = 3897.7646
$ws.Range("M55").Value = -1341.72
$ws.Range("N55").Value = -4243.7646
$ws.Range("H136").Value = 8016.1177
$ws.Range("I136").Value = 7838.6665
$ws.Range("K136").Value = 23515.9995
$ws.Range("M136").Value = -20965.9995
$ws.Range("H140").Value = 82249.55
$ws.Range("J140").Value = 82249.55
$ws.Range("L140").Value = 82249.55
$ws.Range("N140").Value = -92609.55

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 18446.785
$ws.Range("I126").Value = 19675.076
$ws.Range("K126").Value = 59025.228
$ws.Range("M126").Value = -56555.228
$ws.Range("H132").Value = 5978.6895
$ws.Range("I132").Value = 5876.5
$ws.Range("J132").Value = 6299.857
$ws.Range("K132").Value = 17629.5
$ws.Range("L132").Value = 18899.571
$ws.Range("M132").Value = -15099.5
$ws.Range("N132").Value = -23959.571
